$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 is the (previously empty) first data row of the "day sale" table.
# Fill it in with the sold product's details.

# م (serial number) - plain numeric value
$ws.Range("A8").Value = 1

# الاسم (item name) - merged B8:G8
$ws.Range("B8:G8").NumberFormat = "@"
$ws.Range("B8").Value = "HYDRAPHASE UV INTENSE LIGHT CREAM 50 ML"

# الرصيد الحالي (current balance) - merged H8:J8
$ws.Range("H8:J8").NumberFormat = "@"
$ws.Range("H8").Value = "-1:0"

# حد الطلب (order limit) - keep its existing number format, just store
# the quantity as text (matches the source report's export format).
$origFormatK8 = $ws.Range("K8").NumberFormat
$ws.Range("K8").NumberFormat = "@"
$ws.Range("K8").Value = "1"
$ws.Range("K8").NumberFormat = $origFormatK8

# السعر (price)
$ws.Range("L8").NumberFormat = "@"
$ws.Range("L8").Value = "225.00"

# عدد التعاملات (number of transactions)
$ws.Range("M8").NumberFormat = "@"
$ws.Range("M8").Value = "1:0"
